# Logged Week 15 and simulated Week 16
# Updates the Rushing and Receiving stat tables with the new weekly totals,
# then leaves the Receiving sheet active with H10 selected (matching the
# author's last on-screen position).

$wb = $excel.ActiveWorkbook

# ---- Rushing sheet -------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# D.Singletary row (R3) - 1DATT/2DATT/3DATT/RZATT
$rushing.Range("C3").Value = 71
$rushing.Range("D3").Value = 53
$rushing.Range("E3").Value = 2
$rushing.Range("F3").Value = 19

# J.Allen row (R2) - 3DATT
$rushing.Range("E2").Value = 30

# M.Breida row (R5) - 2DATT
$rushing.Range("D5").Value = 11

# I.McKenzie row (R8) - 3DATT
$rushing.Range("E8").Value = 2

# ---- Receiving sheet ------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# D.Singletary row (R2) - Short Target/Short Comp
$receiving.Range("C2").Value = 38
$receiving.Range("D2").Value = 31

# S.Diggs row (R6) - Short Target/Short Comp/Deep Target/Deep Comp/RZ Target/RZ Comp
$receiving.Range("C6").Value = 108
$receiving.Range("D6").Value = 79
$receiving.Range("E6").Value = 33
$receiving.Range("F6").Value = 14
$receiving.Range("G6").Value = 23
$receiving.Range("H6").Value = 14

# E.Sanders row (R7) - Short Target/Short Comp
$receiving.Range("C7").Value = 60
$receiving.Range("D7").Value = 45

# C.Beasley row (R8) - whole row reset to 0 (inactive this week)
$receiving.Range("C8").Value = 0
$receiving.Range("D8").Value = 0
$receiving.Range("E8").Value = 0
$receiving.Range("F8").Value = 0
$receiving.Range("G8").Value = 0
$receiving.Range("H8").Value = 0

# G.Davis row (R9)
$receiving.Range("C9").Value = 50
$receiving.Range("D9").Value = 35
$receiving.Range("E9").Value = 21
$receiving.Range("F9").Value = 15
$receiving.Range("G9").Value = 13
$receiving.Range("H9").Value = 7

# I.McKenzie row (R10) - Short Target/Short Comp
$receiving.Range("C10").Value = 9
$receiving.Range("D10").Value = 7

# D.Knox row (R12) - Short Target/Short Comp/Deep Target
$receiving.Range("C12").Value = 45
$receiving.Range("D12").Value = 37
$receiving.Range("E12").Value = 17

# T.Sweeney row (R13) - Short Target
$receiving.Range("C13").Value = 10

# ---- View state: Receiving tab active, H10 selected -----------------------
$receiving.Activate()
$receiving.Range("H10").Select()
